$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New vm_pu values for rows 2-25 (data rows), columns B-F and I-N.
# Keyed by row number -> ordered array of values for columns B,C,D,E,F,I,J,K,L,M,N
$data = @{
    2 = @(1.02, 1.020206359968576, 1.022057806902428, 1.029458181692308, 1.036301516212148, 1.024866737021063, 1.025404780017896, 1.024893308851434, 1.032272034597252, 1.039095648603761, 1.012550559354305)
    3 = @(1.02, 1.02141517849003, 1.023114353652605, 1.030644089757828, 1.037843058120765, 1.025051908037323, 1.026249332521445, 1.025755757306346, 1.033265075610306, 1.040444822604652, 1.012828668783757)
    4 = @(1.02, 1.02219197056058, 1.023793342091176, 1.031403157238177, 1.038823120909426, 1.025161489473237, 1.026789984231106, 1.026308573244147, 1.033898802735502, 1.041299967126139, 1.013006692245016)
    5 = @(1.02, 1.022517256448988, 1.024077683128152, 1.031720300868926, 1.039231001271457, 1.025205113422967, 1.02701589026047, 1.026539732439602, 1.034163122343347, 1.041655223302059, 1.013081074695027)
    6 = @(1.02, 1.022571798856968, 1.024125360729703, 1.031773435737539, 1.039299244405169, 1.025212294925566, 1.027053740049851, 1.026578472413153, 1.034207380116977, 1.041714624168707, 1.013093537050721)
    7 = @(1.02, 1.022196322047777, 1.023797145799039, 1.031407402638287, 1.038828587244755, 1.025162081973702, 1.026793008220916, 1.026311666880301, 1.033902342814379, 1.041304730724232, 1.013007687943849)
    8 = @(1.02, 1.020616011433081, 1.022415845717053, 1.029860693100594, 1.036826113927952, 1.024931440272417, 1.025691415057484, 1.025185870557923, 1.032609477358089, 1.039555325030591, 1.012644950189477)
    9 = @(1.02, 1.017789318717754, 1.01994549612421, 1.027070826163576, 1.033162546963605, 1.024446268802651, 1.023705040594673, 1.023161350841544, 1.030262809857539, 1.036334408745721, 1.01199077457053)
    10 = @(1.02, 1.015875627551044, 1.018273315277654, 1.025166365554043, 1.030627081488777, 1.024069298437709, 1.02234951634915, 1.021783468249857, 1.028651123295576, 1.034091987348629, 1.011544296912382)
    11 = @(1.02, 1.015039803282574, 1.017543037200454, 1.024330836204852, 1.029506542180896, 1.023893219417235, 1.021754930079198, 1.021179940049249, 1.027941745298043, 1.033097886177078, 1.011348440251087)
    12 = @(1.02, 1.01472823957028, 1.017270826821672, 1.024018819424981, 1.029086865561778, 1.023825870760704, 1.02153290767634, 1.020954708315865, 1.027676494686038, 1.032725109414528, 1.011275304030758)
    13 = @(1.02, 1.014795121251738, 1.017329260298364, 1.024085823858872, 1.02917704500971, 1.02384040554532, 1.02158058535083, 1.021003069312787, 1.02773347176211, 1.032805231571921, 1.011291009587544)
    14 = @(1.02, 1.015014071928757, 1.017520555752961, 1.024305078916408, 1.029471922450175, 1.023887692124988, 1.021736601552811, 1.02116134394565, 1.027919855585727, 1.033067144509454, 1.011342402707767)
    15 = @(1.02, 1.01514882804355, 1.017638292377185, 1.024439947817834, 1.0296531462561, 1.023916568744144, 1.021832573098626, 1.021258721871267, 1.028034459276922, 1.033228049288325, 1.011374016314519)
    16 = @(1.02, 1.015930945245984, 1.01832164891918, 1.025221585141292, 1.030700965982381, 1.024080712292306, 1.022388814577595, 1.021823375587979, 1.028697957648124, 1.034157471292919, 1.011557241455187)
    17 = @(1.02, 1.016419608243189, 1.018748623880118, 1.025708951604482, 1.031352130755798, 1.02418022531991, 1.022735672570994, 1.022175709006006, 1.029111053152437, 1.034734249414734, 1.011671492285435)
    18 = @(1.02, 1.016703945571227, 1.018997073066922, 1.025992175960166, 1.03172976096854, 1.024237030973169, 1.022937253647956, 1.022380555072615, 1.029350896432591, 1.035068447346764, 1.011737889312819)
    19 = @(1.02, 1.016800780609358, 1.019081686938627, 1.026088571184201, 1.031858154367249, 1.024256190569387, 1.023005863435017, 1.022450290195136, 1.029432489599611, 1.035182023940206, 1.011760487861415)
    20 = @(1.02, 1.016367251013211, 1.018702875505945, 1.025656770423228, 1.03128249312955, 1.024169676732061, 1.022698534186213, 1.022137975781349, 1.029066846747035, 1.034672597281755, 1.011659259489206)
    21 = @(1.02, 1.014949626998123, 1.017464250464804, 1.024240559938152, 1.029385184304927, 1.023873821223951, 1.021690691035205, 1.021114765273417, 1.027865018896922, 1.032990115355817, 1.011327279432711)
    22 = @(1.02, 1.014051926558829, 1.016679956109357, 1.023340489607611, 1.0281722329909, 1.023676543338773, 1.02105026008835, 1.020465321391271, 1.027099207490928, 1.031911860295515, 1.011116311746084)
    23 = @(1.02, 1.014528427367476, 1.017096255455992, 1.023818557891494, 1.028817159344181, 1.023782196819635, 1.021390412203231, 1.020810189546404, 1.02750615254289, 1.03248541604529, 1.011228364136729)
    24 = @(1.02, 1.016390911127785, 1.018723549074417, 1.025680352090852, 1.03131396611886, 1.024174447011528, 1.022715317695058, 1.02215502785258, 1.029086825142385, 1.03470046210588, 1.011664787715577)
    25 = @(1.02, 1.018525161123506, 1.02058853013076, 1.027799827053091, 1.034125891266988, 1.024581084450274, 1.024224011147354, 1.023689644860536, 1.034125891266988, 1.037183691109768, 1.012161698658267)
}

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

foreach ($row in $data.Keys) {
    $values = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$row").Value = $values[$i]
    }
}
